$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Status(Summary)")

# ---------------------------------------------------------------------------
# 1) Fix typo in the "To do:" comment of row 6 (D6): "Surevy" -> "Survey".
#    The cell holds rich text (3 runs: plain / bold "To do:" / plain), so we
#    edit just the mis-spelled substring via Characters() and then re-assert
#    the bold run's formatting so the engine keeps the run split instead of
#    collapsing the whole string to a single plain run.
# ---------------------------------------------------------------------------
$d6 = $ws.Range("D6")
$fullText = $d6.Characters().Text
$oldWord = "Surevy"
$newWord = "Survey"
$idx = $fullText.IndexOf($oldWord)
if ($idx -ge 0) {
    $d6.Characters($idx + 1, $oldWord.Length).Text = $newWord
}

$todoLabel = "To do:"
$todoIdx = $fullText.IndexOf($todoLabel)
if ($todoIdx -ge 0) {
    $boldRun = $d6.Characters($todoIdx + 1, $todoLabel.Length)
    $boldRun.Font.Name = "Calibri"
    $boldRun.Font.Size = 11
    $boldRun.Font.Bold = $true
}

$afterTodoStart = $todoIdx + $todoLabel.Length + 1
$afterLen = $fullText.Length - ($afterTodoStart - 1)
if ($afterLen -gt 0) {
    $restRun = $d6.Characters($afterTodoStart, $afterLen)
    $restRun.Font.Name = "Calibri"
    $restRun.Font.Size = 11
    $restRun.Font.Bold = $false
}

# ---------------------------------------------------------------------------
# 2) Populate the previously-empty row 7 with the new status entry.
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = 43902
$ws.Range("B7").Value = "Implementation"
$ws.Range("C7").Value = "SurveyComponent implementation"

$d7 = @'
1) Implemented and tested SurveyComponent outline classes
2) Implemented SurveyResponseItem outline classes completing ouline of all respective class models in the Survey
3) Studied expression-eval, engine.ts and made notes of doubts
4) resolveContent in LocalisedObject expression implemented internally during the object creation
5) SelectionMethod implemented for a dummy map of array items
6) Created an issue for a `wiki` page. Will be updating in due course of `engine` development
'@
$d7 = $d7.TrimEnd("`r", "`n")
$ws.Range("D7").Value = $d7

$ws.Rows.Item(7).RowHeight = 129

# ---------------------------------------------------------------------------
# 3) Scroll the view down and move the active selection to D8.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("D8").Select()
